$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row of data: BanDieuHanh / 12345 / super
$ws.Range("A11").Value = "BanDieuHanh"
$ws.Range("B11").Value = 12345
$ws.Range("C11").Value = "super"

# Update the selection to match the new active cell
$ws.Range("A11").Select()
